$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# "CuadroTexto 2" is the textbox listing library versions (shape 2 on
# this slide); look it up by name so the script is robust to shape
# reordering.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "CuadroTexto 2") {
        $shp = $s.Shapes.Item($i)
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(2)
}

$tr = $shp.TextFrame.TextRange

# Existing text: "Librería pandas versión 2.1.4"
# Add two new paragraphs, each built from three runs so the middle
# (technical / foreign) word stays its own run.
[void]$tr.InsertAfter("`rLibrería ")
[void]$tr.InsertAfter("seaborn")
[void]$tr.InsertAfter(" versión 0.13.2")

[void]$tr.InsertAfter("`rLibrería ")
[void]$tr.InsertAfter("scikit-learn")
[void]$tr.InsertAfter(" versión 1.4.0")
